# Scheduled data refresh: update market-price / profit columns (H:N) across
# the per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the
# latest pulled values. A handful of cells that evaluated to zero-width /
# no-longer-applicable profit figures are cleared outright rather than
# written as 0, matching the upstream scraper's output for this run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 277.3889
$ws.Range("I2").Value = 132.33333
$ws.Range("J2").Value = 349.91666
$ws.Range("K2").Value = 132.33333
$ws.Range("L2").Value = 349.91666
$ws.Range("M2").Value = -19.33332999999999
$ws.Range("N2").Value = -575.91666
$ws.Range("H9").Value = 120
$ws.Range("I9").Value = 120
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 120
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 49
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 390.81818
$ws.Range("I12").Value = 283.33334
$ws.Range("K12").Value = 283.33334
$ws.Range("M12").Value = -113.33334
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H41").Value = 547.94446
$ws.Range("J41").Value = 1072.25
$ws.Range("L41").Value = 1072.25
$ws.Range("N41").Value = -1952.25
$ws.Range("H42").Value = 62.52941
$ws.Range("I42").Value = 63.3125
$ws.Range("K42").Value = 189.9375
$ws.Range("M42").Value = 40.0625
$ws.Range("H43").Value = 12667
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12667
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12667
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12805
$ws.Range("H55").Value = 740.82355
$ws.Range("I55").Value = 392.3125
$ws.Range("J55").Value = 900.1429000000001
$ws.Range("K55").Value = 392.3125
$ws.Range("L55").Value = 900.1429000000001
$ws.Range("M55").Value = -178.3125
$ws.Range("N55").Value = -1328.1429
$ws.Range("H70").Value = 7015.4546
$ws.Range("J70").Value = 8410.25
$ws.Range("L70").Value = 25230.75
$ws.Range("N70").Value = -25770.75
$ws.Range("H73").Value = 7015.4546
$ws.Range("J73").Value = 8410.25
$ws.Range("L73").Value = 25230.75
$ws.Range("N73").Value = -27102.75
$ws.Range("H116").Value = 7591.6113
$ws.Range("I116").Value = 6908.3
$ws.Range("J116").Value = 8445.75
$ws.Range("K116").Value = 6908.3
$ws.Range("L116").Value = 8445.75
$ws.Range("M116").Value = -3466.3
$ws.Range("N116").Value = -15329.75
$ws.Range("H125").Value = 2217.8
$ws.Range("I125").Value = 1894
$ws.Range("K125").Value = 17046
$ws.Range("M125").Value = -14586
$ws.Range("H132").Value = 2221.923
$ws.Range("I132").Value = 2125.8108
$ws.Range("K132").Value = 6377.432400000001
$ws.Range("M132").Value = -3847.432400000001
$ws.Range("H137").Value = 12661375
$ws.Range("I137").Value = 52634596
$ws.Range("J137").Value = 3188.25
$ws.Range("K137").Value = 157903788
$ws.Range("L137").Value = 9564.75
$ws.Range("M137").Value = -157901238
$ws.Range("N137").Value = -14664.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9466.588
$ws.Range("I32").Value = 7873.102
$ws.Range("J32").Value = 48507
$ws.Range("K32").Value = 7873.102
$ws.Range("L32").Value = 48507
$ws.Range("M32").Value = -7586.102
$ws.Range("N32").Value = -49081
$ws.Range("H61").Value = 4983.4634
$ws.Range("I61").Value = 4533.8057
$ws.Range("K61").Value = 4533.8057
$ws.Range("M61").Value = -4321.8057
$ws.Range("H110").Value = 3814.2693
$ws.Range("I110").Value = 2715.0527
$ws.Range("K110").Value = 2715.0527
$ws.Range("M110").Value = -670.0527000000002
$ws.Range("H132").Value = 3724.4424
$ws.Range("I132").Value = 3321.45
$ws.Range("K132").Value = 9964.349999999999
$ws.Range("M132").Value = -7434.349999999999
$ws.Range("H136").Value = 4983.4634
$ws.Range("I136").Value = 4533.8057
$ws.Range("K136").Value = 13601.4171
$ws.Range("M136").Value = -11051.4171

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3627.375
$ws.Range("I20").Value = 3052.5
$ws.Range("J20").Value = 4202.25
$ws.Range("K20").Value = 3052.5
$ws.Range("L20").Value = 4202.25
$ws.Range("M20").Value = -2805.5
$ws.Range("N20").Value = -4696.25
$ws.Range("H86").Value = 6498.8
$ws.Range("I86").Value = 5116.4
$ws.Range("J86").Value = 7881.2
$ws.Range("K86").Value = 5116.4
$ws.Range("L86").Value = 7881.2
$ws.Range("M86").Value = -3993.4
$ws.Range("N86").Value = -10127.2
$ws.Range("H89").Value = 6498.8
$ws.Range("I89").Value = 5116.4
$ws.Range("J89").Value = 7881.2
$ws.Range("K89").Value = 25582
$ws.Range("L89").Value = 39406
$ws.Range("M89").Value = -19966
$ws.Range("N89").Value = -50638
$ws.Range("H105").Value = 15526.625
$ws.Range("I105").Value = 13938.941
$ws.Range("K105").Value = 13938.941
$ws.Range("M105").Value = -12191.941

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3001.6
$ws.Range("I16").Value = 1674.2
$ws.Range("J16").Value = 4329
$ws.Range("K16").Value = 1674.2
$ws.Range("L16").Value = 4329
$ws.Range("M16").Value = -1387.2
$ws.Range("N16").Value = -4903
$ws.Range("H113").Value = 3001.6
$ws.Range("I113").Value = 1674.2
$ws.Range("J113").Value = 4329
$ws.Range("K113").Value = 1674.2
$ws.Range("L113").Value = 4329
$ws.Range("M113").Value = 495.8
$ws.Range("N113").Value = -8669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2278.9412
$ws.Range("I122").Value = 608
$ws.Range("K122").Value = 5472
$ws.Range("M122").Value = -3022

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1986.8611
$ws.Range("I132").Value = 1470.1212
$ws.Range("K132").Value = 4410.363600000001
$ws.Range("M132").Value = -1880.363600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8679.5
$ws.Range("I40").Value = 8388.923000000001
$ws.Range("J40").Value = 9022.909
$ws.Range("K40").Value = 8388.923000000001
$ws.Range("L40").Value = 9022.909
$ws.Range("M40").Value = -8252.923000000001
$ws.Range("N40").Value = -9294.909
$ws.Range("H122").Value = 217492
$ws.Range("I122").Value = 339114.6
$ws.Range("K122").Value = 1017343.8
$ws.Range("M122").Value = -1014893.8
$ws.Range("H132").Value = 3710.8
$ws.Range("I132").Value = 3262.05
$ws.Range("J132").Value = 7300.8
$ws.Range("K132").Value = 9786.150000000001
$ws.Range("L132").Value = 21902.4
$ws.Range("M132").Value = -7256.150000000001
$ws.Range("N132").Value = -26962.4
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("H136").Value = 7167.222
$ws.Range("I136").Value = 4625
$ws.Range("J136").Value = 9201
$ws.Range("K136").Value = 13875
$ws.Range("L136").Value = 27603
$ws.Range("M136").Value = -11325
$ws.Range("N136").Value = -32703

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3252.182
$ws.Range("I132").Value = 2375.2666
$ws.Range("J132").Value = 5131.2856
$ws.Range("K132").Value = 7125.7998
$ws.Range("L132").Value = 15393.8568
$ws.Range("M132").Value = -4595.7998
$ws.Range("N132").Value = -20453.8568
